$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.191952586174011
$ws.Range("B1").Value = 2.260015249252319
$ws.Range("C1").Value = 6.612365245819092
$ws.Range("D1").Value = 2.302503347396851
$ws.Range("E1").Value = 1.188462853431702
